$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.653087377548218
$ws.Range("B1").Value = 1.79763388633728
$ws.Range("C1").Value = 2.039255380630493
$ws.Range("D1").Value = 2.578201532363892
$ws.Range("E1").Value = 1.728454113006592
